# BUG: read_excel failed with empty rows after MultiIndex header (#40649)
# Add a new worksheet "mi_column_empty_rows" that reproduces a MultiIndex
# column header followed by two blank rows before the data rows.

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "mi_column_empty_rows"

# Row 1: first level of the MultiIndex column header ("a", "b")
$newSheet.Range("A1").Value = "a"
$newSheet.Range("B1").Value = "b"

# Row 2: second level of the MultiIndex column header ("A", "B")
$newSheet.Range("A2").Value = "A"
$newSheet.Range("B2").Value = "B"

# Rows 3 and 4 are intentionally left blank to reproduce the bug scenario.

# Rows 5 and 6: the actual data rows
$newSheet.Range("A5").Value = 1
$newSheet.Range("B5").Value = 3
$newSheet.Range("A6").Value = 2
$newSheet.Range("B6").Value = 4
